# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the first data row
# (row 2) of the per-language report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-21 05:02:38"
$zhcn.Range("H2").Value = "2016-03-21 05:02:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-21 05:02:41"
$dede.Range("H2").Value = "2016-03-21 05:03:03"
